$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Record_Transactions_Template")

# C2: change from string "100.50" to the actual numeric value 100.5
$ws.Range("C2").Value = 100.5

# D2: update the note text to include character count suffix
$ws.Range("D2").Value = "Example transaction note REPLACE HERE(28 Character)"
